$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new top data row for 2022-Q4, pushing
#    the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Shift existing data rows 2..7 down to 3..8 (bottom-up so we never clobber
# a row before we've read it).
for ($r = 7; $r -ge 2; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $summary.Range("B$dstRow").Value = $summary.Range("B$srcRow").Value()
    $summary.Range("C$dstRow").Value = $summary.Range("C$srcRow").Value()
    $summary.Range("D$dstRow").Value = $summary.Range("D$srcRow").Value()
}

# Make sure the freshly-occupied row 8 carries the same look (border/font)
# as the row above it, since it was previously empty.
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)

# Re-sequence the running index in column A (0-based) for rows 2..8.
for ($r = 2; $r -le 8; $r++) {
    $summary.Range("A$r").Value = $r - 2
}

# Write the new 2022-Q4 row at the top of the data (row 2).
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" sheet right after "总计" (i.e. right before
#    "2022-Q1"), holding the fund holdings detail for the new quarter.
#    Copying the "2022-Q1" sheet keeps every style/format identical, so we
#    only need to overwrite the actual data afterwards.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q1")
$template.Copy($template)
$newSheet = $wb.Worksheets.Item("2022-Q1 (2)")
$newSheet.Name = "2022-Q4"

# The template sheet only had one data row; make sure we have two (rows 2
# and 3) formatted the same way before writing values into them.
$newSheet.Range("A2:H2").Copy()
$newSheet.Range("A3:H3").PasteSpecial(-4122)

# Row 2 / 3 "index" column stays numeric.
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1

# Columns B..G hold text-formatted values (fund code / name / numbers kept
# as text, matching the rest of the workbook's convention); H is numeric.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("B2").Value = "970073"
$newSheet.Range("C2").Value = "东证融汇成长优选混合A"
$newSheet.Range("D2").Value = "0.38"
$newSheet.Range("E2").Value = "89.59"
$newSheet.Range("F2").Value = "0.84"
$newSheet.Range("G2").Value = "0.0032"
$newSheet.Range("H2").Value = 6

$newSheet.Range("B3").Value = "970074"
$newSheet.Range("C3").Value = "东证融汇成长优选混合C"
$newSheet.Range("D3").Value = "0.11"
$newSheet.Range("E3").Value = "89.59"
$newSheet.Range("F3").Value = "0.84"
$newSheet.Range("G3").Value = "0.0009"
$newSheet.Range("H3").Value = 6

# Restore the originally active tab ("总计"), since copying a sheet shifts
# focus onto the copy.
$wb.Worksheets.Item(1).Activate()

